{"js": "// Replace the 25 \"NN\u00f7N=\" division prompts in the single table with new\n// values, per the commit's regenerated worksheet numbers. Each of the 5\n// populated table rows (rows 0, 4, 8, 12, 16 in the table \u2014 separated by\n// blank spacer rows) holds 5 cells; we target each cell positionally\n// (table row/column index) rather than by text search, because several\n// of the new values collide with old values used elsewhere in the table\n// (e.g. \"50\u00f76=\" is both a pre-edit value in one cell and a post-edit\n// value in another), which would make a simple global search-and-replace\n// ambiguous/unsafe.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"52\u00f75=\", newText: \"62\u00f72=\" },\n  { row: 0, col: 1, oldText: \"96\u00f73=\", newText: \"75\u00f76=\" },\n  { row: 0, col: 2, oldText: \"65\u00f74=\", newText: \"16\u00f78=\" },\n  { row: 0, col: 3, oldText: \"75\u00f73=\", newText: \"24\u00f72=\" },\n  { row: 0, col: 4, oldText: \"46\u00f72=\", newText: \"64\u00f75=\" },\n  { row: 1, col: 0, oldText: \"96\u00f74=\", newText: \"83\u00f77=\" },\n  { row: 1, col: 1, oldText: \"13\u00f77=\", newText: \"56\u00f78=\" },\n  { row: 1, col: 2, oldText: \"94\u00f76=\", newText: \"23\u00f76=\" },\n  { row: 1, col: 3, oldText: \"18\u00f78=\", newText: \"18\u00f77=\" },\n  { row: 1, col: 4, oldText: \"51\u00f76=\", newText: \"71\u00f74=\" },\n  { row: 2, col: 0, oldText: \"57\u00f73=\", newText: \"50\u00f76=\" },\n  { row: 2, col: 1, oldText: \"95\u00f78=\", newText: \"85\u00f77=\" },\n  { row: 2, col: 2, oldText: \"15\u00f75=\", newText: \"44\u00f77=\" },\n  { row: 2, col: 3, oldText: \"41\u00f77=\", newText: \"88\u00f74=\" },\n  { row: 2, col: 4, oldText: \"10\u00f72=\", newText: \"51\u00f75=\" },\n  { row: 3, col: 0, oldText: \"20\u00f79=\", newText: \"77\u00f77=\" },\n  { row: 3, col: 1, oldText: \"86\u00f73=\", newText: \"98\u00f77=\" },\n  { row: 3, col: 2, oldText: \"73\u00f72=\", newText: \"99\u00f77=\" },\n  { row: 3, col: 3, oldText: \"89\u00f73=\", newText: \"28\u00f72=\" },\n  { row: 3, col: 4, oldText: \"83\u00f73=\", newText: \"95\u00f78=\" },\n  { row: 4, col: 0, oldText: \"37\u00f79=\", newText: \"86\u00f76=\" },\n  { row: 4, col: 1, oldText: \"50\u00f76=\", newText: \"13\u00f76=\" },\n  { row: 4, col: 2, oldText: \"33\u00f77=\", newText: \"97\u00f74=\" },\n  { row: 4, col: 3, oldText: \"57\u00f72=\", newText: \"84\u00f72=\" },\n  { row: 4, col: 4, oldText: \"71\u00f78=\", newText: \"93\u00f76=\" },\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected worksheet table not found in document body.\");\n}\n\n// Data rows sit every 4th table row (1 data row + 3 blank spacer rows).\nconst ROW_STRIDE = 4;\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row * ROW_STRIDE, col);\n  const cellRange = cell.body.getRange();\n  cellRange.load(\"text\");\n  await context.sync();\n\n  const currentText = cellRange.text.trim();\n  if (currentText !== oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${row}, col ${col}: expected \"${oldText}\" but found \"${currentText}\"`\n    );\n  }\n\n  cellRange.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"NN\u00f7N=\" division prompts in the single table with new\n# values, per the commit's regenerated worksheet numbers. The table has\n# 5 populated rows (table rows 1, 5, 9, 13, 17 \u2014 1-indexed, separated by\n# 3 blank spacer rows each) with 5 columns apiece. We target each cell\n# positionally via Table.Cell(row, col) rather than Find/Replace on text,\n# because several of the new values collide with old values used\n# elsewhere in the table (e.g. \"50\u00f76=\" is both a pre-edit value in one\n# cell and a post-edit value in another), which would make a simple\n# global find-and-replace ambiguous/unsafe.\n\n$d = $word.ActiveDocument\n$table = $d.Tables(1)\n\n$replacements = @(\n  @{ Row = 1;  Col = 1; Old = \"52\u00f75=\"; New = \"62\u00f72=\" },\n  @{ Row = 1;  Col = 2; Old = \"96\u00f73=\"; New = \"75\u00f76=\" },\n  @{ Row = 1;  Col = 3; Old = \"65\u00f74=\"; New = \"16\u00f78=\" },\n  @{ Row = 1;  Col = 4; Old = \"75\u00f73=\"; New = \"24\u00f72=\" },\n  @{ Row = 1;  Col = 5; Old = \"46\u00f72=\"; New = \"64\u00f75=\" },\n  @{ Row = 5;  Col = 1; Old = \"96\u00f74=\"; New = \"83\u00f77=\" },\n  @{ Row = 5;  Col = 2; Old = \"13\u00f77=\"; New = \"56\u00f78=\" },\n  @{ Row = 5;  Col = 3; Old = \"94\u00f76=\"; New = \"23\u00f76=\" },\n  @{ Row = 5;  Col = 4; Old = \"18\u00f78=\"; New = \"18\u00f77=\" },\n  @{ Row = 5;  Col = 5; Old = \"51\u00f76=\"; New = \"71\u00f74=\" },\n  @{ Row = 9;  Col = 1; Old = \"57\u00f73=\"; New = \"50\u00f76=\" },\n  @{ Row = 9;  Col = 2; Old = \"95\u00f78=\"; New = \"85\u00f77=\" },\n  @{ Row = 9;  Col = 3; Old = \"15\u00f75=\"; New = \"44\u00f77=\" },\n  @{ Row = 9;  Col = 4; Old = \"41\u00f77=\"; New = \"88\u00f74=\" },\n  @{ Row = 9;  Col = 5; Old = \"10\u00f72=\"; New = \"51\u00f75=\" },\n  @{ Row = 13; Col = 1; Old = \"20\u00f79=\"; New = \"77\u00f77=\" },\n  @{ Row = 13; Col = 2; Old = \"86\u00f73=\"; New = \"98\u00f77=\" },\n  @{ Row = 13; Col = 3; Old = \"73\u00f72=\"; New = \"99\u00f77=\" },\n  @{ Row = 13; Col = 4; Old = \"89\u00f73=\"; New = \"28\u00f72=\" },\n  @{ Row = 13; Col = 5; Old = \"83\u00f73=\"; New = \"95\u00f78=\" },\n  @{ Row = 17; Col = 1; Old = \"37\u00f79=\"; New = \"86\u00f76=\" },\n  @{ Row = 17; Col = 2; Old = \"50\u00f76=\"; New = \"13\u00f76=\" },\n  @{ Row = 17; Col = 3; Old = \"33\u00f77=\"; New = \"97\u00f74=\" },\n  @{ Row = 17; Col = 4; Old = \"57\u00f72=\"; New = \"84\u00f72=\" },\n  @{ Row = 17; Col = 5; Old = \"71\u00f78=\"; New = \"93\u00f76=\" }\n)\n\nforeach ($r in $replacements) {\n  $cell = $table.Cell($r.Row, $r.Col)\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $r.Old) {\n    throw \"Unexpected cell text at row $($r.Row), col $($r.Col): expected '$($r.Old)' but found '$current'\"\n  }\n  $cell.Range.Text = $r.New\n}\n"}
